$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("D12").Value = 0
$ws.Range("E12").Formula = "=1000*0.000678331591188908"

$plotSheet = $wb.Worksheets.Item("Dynamic Power Plots")
$co = $plotSheet.ChartObjects(1)
$chart = $co.Chart
$series = $chart.FullSeriesCollection(2)
Write-Host "before:" $series.Formula
try {
  $series.Formula = "=SERIES(Data!`$E`$3,Data!`$B`$5:`$B`$23,Data!`$E`$5:`$E`$23,2)"
  Write-Host "set formula ok"
} catch {
  Write-Host "set formula failed: $_"
}
Write-Host "after:" $series.Formula
Write-Host $series.Values
